$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.554.32"
$ws.Range("E2").Value = "  -0.47%  "
$ws.Range("D3").Value = "2.301.80"
$ws.Range("E3").Value = "  +0.22%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "537.06"
$ws.Range("E5").Value = "  -2.05%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "132.05"
$ws.Range("E6").Value = "  +0.89%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.585"
$ws.Range("E8").Value = "  +2.30%  "
$ws.Range("D9").Value = "2.300.37"
$ws.Range("E9").Value = "  +0.26%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.100"
$ws.Range("E10").Value = "  -1.69%  "
$ws.Range("E11").Value = "  -1.12%  "
$ws.Range("E12").Value = "  +0.70%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.333"
$ws.Range("E13").Value = "  -0.85%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "23.69"
$ws.Range("E14").Value = "  -1.03%  "
$ws.Range("D15").Value = "2.712.75"
$ws.Range("E15").Value = "  +0.35%  "
$ws.Range("D16").Value = "58.471.62"
$ws.Range("E16").Value = "  -0.53%  "
$ws.Range("E17").Value = "  -0.65%  "
$ws.Range("D18").Value = "2.300.37"
$ws.Range("E18").Value = "  +3.13%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.56"
$ws.Range("E19").Value = "  -1.15%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.19"
$ws.Range("E20").Value = "  -2.81%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "315.76"
$ws.Range("E21").Value = "  +0.12%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.58"
$ws.Range("E22").Value = "  +1.50%  "
$ws.Range("E23").Value = "  +0.00%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "63.03"
$ws.Range("E24").Value = "  -0.22%  "
$ws.Range("E25").Value = "  -1.93%  "
$ws.Range("E26").Value = "  +0.14%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.94"
$ws.Range("E27").Value = "  -2.06%  "
$ws.Range("E28").Value = "  -0.80%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "171.51"
$ws.Range("E29").Value = "  +1.08%  "
$ws.Range("E30").Value = "  -2.53%  "
$ws.Range("D31").Value = "0.0₃0725"
$ws.Range("E31").Value = "  -0.66%  "
$ws.Range("E32").Value = "  -0.78%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.83"
$ws.Range("E33").Value = "  +0.45%  "
$ws.Range("E34").Value = "  -0.08%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "17.89"
$ws.Range("E36").Value = "  +0.66%  "
$ws.Range("E37").Value = "  -0.03%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.25"
$ws.Range("E38").Value = "  -0.57%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.99"
$ws.Range("E39").Value = "  +0.48%  "
$ws.Range("E40").Value = "  -0.64%  "
$ws.Range("B41").Value = "Aave"
$ws.Range("C41").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "141.24"
$ws.Range("E41").Value = "  +0.85%  "
$ws.Range("B42").Value = "Bittensor"
$ws.Range("C42").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "289.72"
$ws.Range("E42").Value = "  -4.10%  "
$ws.Range("E43").Value = "  -0.05%  "
$ws.Range("E44").Value = "  -0.11%  "
$ws.Range("E45").Value = "  -1.14%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.556"
$ws.Range("E46").Value = "  -0.23%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "18.28"
$ws.Range("E47").Value = "  -2.00%  "
$ws.Range("E48").Value = "  -2.79%  "
$ws.Range("E49").Value = "  -0.60%  "
$ws.Range("E50").Value = "  +0.12%  "
$ws.Range("E51").Value = "  -0.30%  "

Write-Output "Applied cryptos update"
